# Update cryptocurrency Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.339.55"
$ws.Range("E2").Value = "  +0.96%  "

$ws.Range("D3").Value = "3.355.90"
$ws.Range("E3").Value = "  +0.71%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'585.11"
$ws.Range("E5").Value = "  +0.82%  "

$ws.Range("D6").Value = "'177.84"
$ws.Range("E6").Value = "  +1.24%  "

$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("E8").Value = "  +0.34%  "

$ws.Range("D9").Value = "'0.185"
$ws.Range("E9").Value = "  +3.82%  "

$ws.Range("D10").Value = "'0.581"
$ws.Range("E10").Value = "  +0.89%  "

$ws.Range("E11").Value = "  +5.85%  "

$ws.Range("E12").Value = "  +1.63%  "

$ws.Range("D13").Value = "'689.24"
$ws.Range("E13").Value = "  +2.90%  "

$ws.Range("D14").Value = "3.906.38"
$ws.Range("E14").Value = "  +0.73%  "

$ws.Range("D16").Value = "68.373.38"
$ws.Range("E16").Value = "  +1.04%  "

$ws.Range("E17").Value = "  +1.37%  "

$ws.Range("D18").Value = "3.386.76"
$ws.Range("E18").Value = "  +1.57%  "

$ws.Range("E19").Value = "  +0.62%  "

$ws.Range("E20").Value = "  +2.35%  "

$ws.Range("E21").Value = "  +0.79%  "

$ws.Range("D22").Value = "'5.47"
$ws.Range("E22").Value = "  +0.98%  "

$ws.Range("E23").Value = "  -0.75%  "

$ws.Range("D24").Value = "'100.07"
$ws.Range("E24").Value = "  +1.27%  "

$ws.Range("E25").Value = "  +1.72%  "

$ws.Range("E26").Value = "  +1.25%  "

$ws.Range("E27").Value = "  +2.83%  "

$ws.Range("D28").Value = "'32.98"
$ws.Range("E28").Value = "  -1.77%  "

$ws.Range("E29").Value = "  +1.36%  "

$ws.Range("D30").Value = "'6.96"
$ws.Range("E30").Value = "  -5.01%  "

$ws.Range("D31").Value = "'11.08"
$ws.Range("E31").Value = "  +1.20%  "

$ws.Range("D32").Value = "'552.99"
$ws.Range("E32").Value = "  -3.52%  "

$ws.Range("D33").Value = "'0.105"
$ws.Range("E33").Value = "  +0.76%  "

$ws.Range("D34").Value = "'58.07"
$ws.Range("E34").Value = "  +2.56%  "

$ws.Range("D35").Value = "3.719.06"
$ws.Range("E35").Value = "  +0.90%  "

$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("E37").Value = "  +1.33%  "

$ws.Range("E38").Value = "  +8.21%  "

$ws.Range("D39").Value = "'34.75"
$ws.Range("E39").Value = "  +1.19%  "

$ws.Range("E40").Value = "  +2.26%  "

$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("E42").Value = "  +1.33%  "

$ws.Range("E43").Value = "  +0.72%  "

$ws.Range("D44").Value = "'3.28"
$ws.Range("E44").Value = "  -1.38%  "

$ws.Range("E45").Value = "  +1.55%  "

$ws.Range("E46").Value = "  +1.92%  "

$ws.Range("E47").Value = "  +0.43%  "

$ws.Range("E48").Value = "  -0.01%  "

$ws.Range("E49").Value = "  -0.42%  "

$ws.Range("D50").Value = "'131.72"
$ws.Range("E50").Value = "  +2.11%  "

$ws.Range("E51").Value = "  -1.25%  "

